# Refresh the crypto price/volume table with the latest scraped values.
# (GitHub Actions data-refresh run.)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.181.48'
$ws.Range('E2').Value = '  +0.38%  '
$ws.Range('D3').Value = '1.835.91'
$ws.Range('E3').Value = '  +0.61%  '
$cell = $ws.Range('D4')
$cell.NumberFormat = '@'
$cell.Value = '1.007'
$cell.Style = 'Normal'
$ws.Range('E4').Value = '  -0.08%  '
$cell = $ws.Range('D5')
$cell.NumberFormat = '@'
$cell.Value = '312.88'
$cell.Style = 'Normal'
$ws.Range('E5').Value = '  +0.19%  '
$cell = $ws.Range('D6')
$cell.NumberFormat = '@'
$cell.Value = '1.006'
$cell.Style = 'Normal'
$ws.Range('E6').Value = '  -0.04%  '
$cell = $ws.Range('D7')
$cell.NumberFormat = '@'
$cell.Value = '0.4647'
$cell.Style = 'Normal'
$ws.Range('E7').Value = '  -0.95%  '
$cell = $ws.Range('D8')
$cell.NumberFormat = '@'
$cell.Value = '0.3718'
$cell.Style = 'Normal'
$ws.Range('E8').Value = '  +1.88%  '
$cell = $ws.Range('D9')
$cell.NumberFormat = '@'
$cell.Value = '0.07367'
$cell.Style = 'Normal'
$ws.Range('E9').Value = '  -0.30%  '
$cell = $ws.Range('D10')
$cell.NumberFormat = '@'
$cell.Value = '0.8768'
$cell.Style = 'Normal'
$ws.Range('E10').Value = '  -0.17%  '
$ws.Range('E11').Value = '  -1.03%  '
$ws.Range('B12').Value = 'TRON'
$ws.Range('C12').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$cell = $ws.Range('D12')
$cell.NumberFormat = '@'
$cell.Value = '0.07816'
$cell.Style = 'Normal'
$ws.Range('E12').Value = '  +3.66%  '
$ws.Range('B13').Value = 'Chainlink'
$ws.Range('C13').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$cell = $ws.Range('D13')
$cell.NumberFormat = '@'
$cell.Value = '6.623'
$cell.Style = 'Normal'
$ws.Range('E13').Value = '  +1.50%  '
$ws.Range('B14').Value = 'WrappedEther'
$ws.Range('C14').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D14').Value = '1.743.57'
$ws.Range('E14').Value = '  -7.20%  '
$ws.Range('B15').Value = 'Polkadot'
$ws.Range('C15').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$cell = $ws.Range('D15')
$cell.NumberFormat = '@'
$cell.Value = '5.355'
$cell.Style = 'Normal'
$ws.Range('E15').Value = '  -0.35%  '
$cell = $ws.Range('D16')
$cell.NumberFormat = '@'
$cell.Value = '92.18'
$cell.Style = 'Normal'
$ws.Range('E16').Value = '  -0.65%  '
$cell = $ws.Range('D17')
$cell.NumberFormat = '@'
$cell.Value = '1.008'
$cell.Style = 'Normal'
$ws.Range('E17').Value = '  +0.17%  '
$cell = $ws.Range('D18')
$cell.NumberFormat = '@'
$cell.Value = '0.000008862'
$cell.Style = 'Normal'
$ws.Range('E18').Value = '  +1.64%  '
$cell = $ws.Range('D19')
$cell.NumberFormat = '@'
$cell.Value = '1.007'
$cell.Style = 'Normal'
$ws.Range('E19').Value = '  +0.00%  '
$ws.Range('B20').Value = 'WrappedBTC'
$ws.Range('C20').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D20').Value = '27.446.23'
$ws.Range('E20').Value = '  -0.30%  '
$ws.Range('B21').Value = 'Avalanche'
$ws.Range('C21').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$cell = $ws.Range('D21')
$cell.NumberFormat = '@'
$cell.Value = '14.72'
$cell.Style = 'Normal'
$ws.Range('E21').Value = '  +0.78%  '
$cell = $ws.Range('D22')
$cell.NumberFormat = '@'
$cell.Value = '5.152'
$cell.Style = 'Normal'
$cell = $ws.Range('D23')
$cell.NumberFormat = '@'
$cell.Value = '10.61'
$cell.Style = 'Normal'
$ws.Range('E23').Value = '  +0.04%  '
$ws.Range('D24').Value = '2.131.74'
$ws.Range('E24').Value = '  +2.59%  '
$cell = $ws.Range('D25')
$cell.NumberFormat = '@'
$cell.Value = '152.35'
$cell.Style = 'Normal'
$cell = $ws.Range('D26')
$cell.NumberFormat = '@'
$cell.Value = '1.832'
$cell.Style = 'Normal'
$ws.Range('E26').Value = '  -2.34%  '
$cell = $ws.Range('D27')
$cell.NumberFormat = '@'
$cell.Value = '18.42'
$cell.Style = 'Normal'
$ws.Range('E27').Value = '  -0.47%  '
$cell = $ws.Range('D28')
$cell.NumberFormat = '@'
$cell.Value = '2.109'
$cell.Style = 'Normal'
$ws.Range('E28').Value = '  -1.11%  '
$cell = $ws.Range('D29')
$cell.NumberFormat = '@'
$cell.Value = '5.097'
$cell.Style = 'Normal'
$ws.Range('E29').Value = '  -1.34%  '
$cell = $ws.Range('D30')
$cell.NumberFormat = '@'
$cell.Value = '115.91'
$cell.Style = 'Normal'
$ws.Range('E30').Value = '  -0.44%  '
$cell = $ws.Range('D31')
$cell.NumberFormat = '@'
$cell.Value = '0.08882'
$cell.Style = 'Normal'
$ws.Range('E31').Value = '  -0.33%  '
$cell = $ws.Range('D32')
$cell.NumberFormat = '@'
$cell.Value = '2.958'
$cell.Style = 'Normal'
$ws.Range('E32').Value = '  +0.54%  '
$cell = $ws.Range('D33')
$cell.NumberFormat = '@'
$cell.Value = '0.7331'
$cell.Style = 'Normal'
$ws.Range('E33').Value = '  -1.34%  '
$cell = $ws.Range('D34')
$cell.NumberFormat = '@'
$cell.Value = '4.460'
$cell.Style = 'Normal'
$ws.Range('E34').Value = '  -1.08%  '
$cell = $ws.Range('D35')
$cell.NumberFormat = '@'
$cell.Value = '1.144'
$cell.Style = 'Normal'
$ws.Range('E35').Value = '  -1.52%  '
$ws.Range('E36').Value = '  -1.32%  '
$ws.Range('B37').Value = 'TrustWalletToken'
$ws.Range('C37').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$cell = $ws.Range('D37')
$cell.NumberFormat = '@'
$cell.Value = '1.074'
$cell.Style = 'Normal'
$ws.Range('E37').Value = '  -1.44%  '
$ws.Range('B38').Value = 'VeChain'
$ws.Range('C38').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$cell = $ws.Range('D38')
$cell.NumberFormat = '@'
$cell.Value = '0.01952'
$cell.Style = 'Normal'
$ws.Range('E38').Value = '  +1.15%  '
$cell = $ws.Range('D39')
$cell.NumberFormat = '@'
$cell.Value = '0.05241'
$cell.Style = 'Normal'
$ws.Range('E39').Value = '  -1.08%  '
$cell = $ws.Range('D40')
$cell.NumberFormat = '@'
$cell.Value = '2.934'
$cell.Style = 'Normal'
$ws.Range('E40').Value = '  +0.11%  '
$cell = $ws.Range('D41')
$cell.NumberFormat = '@'
$cell.Value = '7.196'
$cell.Style = 'Normal'
$ws.Range('E41').Value = '  -1.65%  '
$cell = $ws.Range('D42')
$cell.NumberFormat = '@'
$cell.Value = '0.5212'
$cell.Style = 'Normal'
$ws.Range('E42').Value = '  -0.64%  '
$cell = $ws.Range('D43')
$cell.NumberFormat = '@'
$cell.Value = '0.8835'
$cell.Style = 'Normal'
$ws.Range('E43').Value = '  -12.28%  '
$cell = $ws.Range('D44')
$cell.NumberFormat = '@'
$cell.Value = '0.1634'
$cell.Style = 'Normal'
$ws.Range('E44').Value = '  -0.36%  '
$cell = $ws.Range('D45')
$cell.NumberFormat = '@'
$cell.Value = '8.256'
$cell.Style = 'Normal'
$ws.Range('E45').Value = '  -1.29%  '
$cell = $ws.Range('D46')
$cell.NumberFormat = '@'
$cell.Value = '0.4859'
$cell.Style = 'Normal'
$ws.Range('E46').Value = '  -0.65%  '
$ws.Range('B47').Value = 'PaxDollar'
$ws.Range('C47').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$cell = $ws.Range('D47')
$cell.NumberFormat = '@'
$cell.Value = '1.007'
$cell.Style = 'Normal'
$ws.Range('E47').Value = '  +0.00%  '
$ws.Range('B48').Value = 'EnergySwap'
$ws.Range('C48').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$cell = $ws.Range('D48')
$cell.NumberFormat = '@'
$cell.Value = '10.20'
$cell.Style = 'Normal'
$ws.Range('E48').Value = '  -2.28%  '
$cell = $ws.Range('D49')
$cell.NumberFormat = '@'
$cell.Value = '102.91'
$cell.Style = 'Normal'
$ws.Range('E49').Value = '  -1.28%  '
$cell = $ws.Range('D50')
$cell.NumberFormat = '@'
$cell.Value = '1.632'
$cell.Style = 'Normal'
$ws.Range('E50').Value = '  -1.06%  '
$cell = $ws.Range('D51')
$cell.NumberFormat = '@'
$cell.Value = '0.06235'
$cell.Style = 'Normal'
$ws.Range('E51').Value = '  -0.53%  '
